$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table on the active sheet lists historical-distance records
# (title / timestamp / historical distance / time bucket / uri).
# A new JSON source record was folded into the time-bucket analysis
# and the table was regenerated, which re-sorted the rows: the
# "News: Taken by Storm" record (previously row 4) now sorts to the
# very end of the table (row 8), and every row below it shifts up by
# one. Re-apply that same row reshuffle here, including moving its
# hyperlink so it keeps following the correct URI.

# Snapshot the 7 data rows (A:E) in their current (pre-edit) order.
$titles = @(
    "The ice storm : an historic record in photographs of January 1998",
    "National Weather Service Forecast Office",
    "News: Taken by Storm",
    "The St. Lawrence River Valley 1998 ice storm: maps and facts",
    "Operations : Past Operations : Operation Recuperation",
    "Canadian agriculture at a glance 1999: Article",
    "Ice Storm 1998"
)
$timestamps = @(
    "1998-01-08T00:00:00UTC",
    "1-01-01T00:00:00UTC",
    "1-01-01T00:00:00UTC",
    "1-01-01T00:00:00UTC",
    "1-01-01T00:00:00UTC",
    "1-01-01T00:00:00UTC",
    "1-01-01T00:00:00UTC"
)
$dist = @(7, "unknown", "unknown", "unknown", "unknown", "unknown", "unknown")
$bucket = @("day_2_to_30", "unknown", "unknown", "unknown", "unknown", "unknown", "unknown")
$uris = @(
    "https://archive.org/details/icestormhistoric0000able",
    "https://web.archive.org/web/20080511204226/http://www.erh.noaa.gov/btv/events/IceStorm1998/ice98.shtml",
    "https://web.archive.org/web/20090612022908/http://www.theweathernetwork.com/news/storm_watch_stories3&stormfile=topstorms2_01_06_2009",
    "http://www.statcan.gc.ca/pub/16f0021x/16f0021x1998001-eng.htm",
    "https://web.archive.org/web/20060529011947/http://www.forces.gc.ca/site/Operations/recuperation_e.asp",
    "https://web.archive.org/web/20060308012756/http://www.statcan.ca/english/kits/agric/ice.htm",
    "http://www.msc-smc.ec.gc.ca/media/icestorm98/icestorm98_the_worst_e.cfm"
)

# New row order: move index 2 ("News: Taken by Storm", old row 4) to
# the end; everything else keeps its relative order.
$order = @(0, 1, 3, 4, 5, 6, 2)

# Clear the existing hyperlinks on the sheet before rewriting the
# uri column so stale relationships don't linger.
$ws.Range("E2").Hyperlinks.Delete()

for ($i = 0; $i -lt $order.Length; $i++) {
    $src = $order[$i]
    $r = 2 + $i

    $ws.Cells.Item($r, 1).Value = $titles[$src]
    $ws.Cells.Item($r, 2).Value = $timestamps[$src]
    $ws.Cells.Item($r, 3).Value = $dist[$src]
    $ws.Cells.Item($r, 4).Value = $bucket[$src]
    $ws.Cells.Item($r, 5).Value = $uris[$src]
}

# Re-create the hyperlinks for column E, top to bottom, so they line
# up with the reshuffled uri text. Re-apply the shared "Hyperlink"
# cell style afterwards so every link cell keeps using the workbook's
# existing named style instead of Hyperlinks.Add's ad-hoc copy.
for ($i = 0; $i -lt $order.Length; $i++) {
    $r = 2 + $i
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 5), $ws.Cells.Item($r, 5).Value2) | Out-Null
    $ws.Cells.Item($r, 5).Style = "Hyperlink"
}
